# Updated cryptos list values (Price and Volume(1h)) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.105.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.225.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '291.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.41%  '

$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0782'
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = '  +3.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.571.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("E15").Value = '  -1.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.214.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.729'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.050.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0888'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.94%  '

$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("E28").Value = '  -2.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '156.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.42%  '

$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("E33").Value = '  +1.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0720'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.87%  '

$ws.Range("E35").Value = '  +7.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.33%  '

$ws.Range("E37").Value = '  +0.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0984'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.128.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.01%  '

$ws.Range("E42").Value = '  +2.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.14%  '

$ws.Range("E44").Value = '  -3.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0268'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("E46").Value = '  +0.64%  '

$ws.Range("E47").Value = '  +3.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.437.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.46%  '
